# Weekly update: prepend two new price records (for Fecha serial 44588) at
# the top of the data block, shifting all existing records (rows 332-427)
# down by two rows to 334-429. This matches the commit "Fruta / hortaliza,
# semanal" and the dimension growing from A1:R427 to A1:R429.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 332, pushing the rest
# of the table (332:427) down to (334:429).
$ws.Rows("332:333").Insert()

# --- New row 332 ---------------------------------------------------------
$ws.Cells.Item(332, 1).Value  = 4
$ws.Cells.Item(332, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(332, 3).Value  = "Los Lagos"
$ws.Cells.Item(332, 4).Value  = 44588
$ws.Cells.Item(332, 5).Value  = 10
$ws.Cells.Item(332, 6).Value  = 100112004
$ws.Cells.Item(332, 7).Value  = "Cebolla"
$ws.Cells.Item(332, 8).Value  = "Morada(o)"
$ws.Cells.Item(332, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(332, 10).Value = 90
$ws.Cells.Item(332, 11).Value = 15000
$ws.Cells.Item(332, 12).Value = 15000
$ws.Cells.Item(332, 13).Value = 15000
$ws.Cells.Item(332, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(332, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(332, 16).Value = 833
$ws.Cells.Item(332, 17).Value = 18
$ws.Cells.Item(332, 18).Value = "Hortaliza"

# --- New row 333 ---------------------------------------------------------
$ws.Cells.Item(333, 1).Value  = 4
$ws.Cells.Item(333, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(333, 3).Value  = "Los Lagos"
$ws.Cells.Item(333, 4).Value  = 44588
$ws.Cells.Item(333, 5).Value  = 10
$ws.Cells.Item(333, 6).Value  = 100112004
$ws.Cells.Item(333, 7).Value  = "Cebolla"
$ws.Cells.Item(333, 8).Value  = "Sin especificar"
$ws.Cells.Item(333, 9).Value  = "1a (cosecha)"
$ws.Cells.Item(333, 10).Value = 300
$ws.Cells.Item(333, 11).Value = 6500
$ws.Cells.Item(333, 12).Value = 6500
$ws.Cells.Item(333, 13).Value = 6500
$ws.Cells.Item(333, 14).Value = "$/malla 16 kilos"
$ws.Cells.Item(333, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(333, 16).Value = 406
$ws.Cells.Item(333, 17).Value = 16
$ws.Cells.Item(333, 18).Value = "Hortaliza"
